$d = $word.ActiveDocument

$found1 = $d.Content.Find.Execute("extensionistas exitosas.Elaboração, desenvolvimento", $true, $false, $false, $false, $false, $true, 1, $false, "extensionistas exitosas.^lElaboração, desenvolvimento", 2)
if (-not $found1) { Write-Host "NOT FOUND #1: extensionistas exitosas.Elaboração, desenvolvimento" }
$found2 = $d.Content.Find.Execute("ampliação do conhecimento.Aulas expositivas interativas online", $true, $false, $false, $false, $false, $true, 1, $false, "ampliação do conhecimento.^lAulas expositivas interativas online", 2)
if (-not $found2) { Write-Host "NOT FOUND #2: ampliação do conhecimento.Aulas expositivas interativas online" }
$found3 = $d.Content.Find.Execute("no início dos projetos.[1] CASARIN", $true, $false, $false, $false, $false, $true, 1, $false, "no início dos projetos.^l^l[1] CASARIN", 2)
if (-not $found3) { Write-Host "NOT FOUND #3: no início dos projetos.[1] CASARIN" }
$found4 = $d.Content.Find.Execute("InterSaberes, 2012.[2] CORDEIRO", $true, $false, $false, $false, $false, $true, 1, $false, "InterSaberes, 2012.^l[2] CORDEIRO", 2)
if (-not $found4) { Write-Host "NOT FOUND #4: InterSaberes, 2012.[2] CORDEIRO" }
$found5 = $d.Content.Find.Execute("InterSaberes, 2012. [3] PEROVANO", $true, $false, $false, $false, $false, $true, 1, $false, "InterSaberes, 2012. ^l[3] PEROVANO", 2)
if (-not $found5) { Write-Host "NOT FOUND #5: InterSaberes, 2012. [3] PEROVANO" }
$found6 = $d.Content.Find.Execute("InterSaberes, 2016. [4] BARROS", $true, $false, $false, $false, $false, $true, 1, $false, "InterSaberes, 2016. ^l[4] BARROS", 2)
if (-not $found6) { Write-Host "NOT FOUND #6: InterSaberes, 2016. [4] BARROS" }
$found7 = $d.Content.Find.Execute("Hall, 2007.[5] CASTRO", $true, $false, $false, $false, $false, $true, 1, $false, "Hall, 2007.^l[5] CASTRO", 2)
if (-not $found7) { Write-Host "NOT FOUND #7: Hall, 2007.[5] CASTRO" }
$found8 = $d.Content.Find.Execute("Hall, 2011.[6] FOGGETTI", $true, $false, $false, $false, $false, $true, 1, $false, "Hall, 2011.^l[6] FOGGETTI", 2)
if (-not $found8) { Write-Host "NOT FOUND #8: Hall, 2011.[6] FOGGETTI" }
$found9 = $d.Content.Find.Execute("Pearson, 2014. [7] MAGALHÃES", $true, $false, $false, $false, $false, $true, 1, $false, "Pearson, 2014. ^l[7] MAGALHÃES", 2)
if (-not $found9) { Write-Host "NOT FOUND #9: Pearson, 2014. [7] MAGALHÃES" }
$found10 = $d.Content.Find.Execute("Ática, 2005. [8] MATTAR", $true, $false, $false, $false, $false, $true, 1, $false, "Ática, 2005. ^l[8] MATTAR", 2)
if (-not $found10) { Write-Host "NOT FOUND #10: Ática, 2005. [8] MATTAR" }
$found11 = $d.Content.Find.Execute("Saraiva, 2013. [9] BRASIL", $true, $false, $false, $false, $false, $true, 1, $false, "Saraiva, 2013. ^l[9] BRASIL", 2)
if (-not $found11) { Write-Host "NOT FOUND #11: Saraiva, 2013. [9] BRASIL" }
$found12 = $d.Content.Find.Execute("ISBN: 9788585819071.[11] GIL", $true, $false, $false, $false, $false, $true, 1, $false, "ISBN: 9788585819071.^l[11] GIL", 2)
if (-not $found12) { Write-Host "NOT FOUND #12: ISBN: 9788585819071.[11] GIL" }
$found13 = $d.Content.Find.Execute("9788522458233.[12] RAMALHO", $true, $false, $false, $false, $false, $true, 1, $false, "9788522458233.^l[12] RAMALHO", 2)
if (-not $found13) { Write-Host "NOT FOUND #13: 9788522458233.[12] RAMALHO" }
$found14 = $d.Content.Find.Execute("formato de arquivo.[13] ABREU", $true, $false, $false, $false, $false, $true, 1, $false, "formato de arquivo.^l[13] ABREU", 2)
if (-not $found14) { Write-Host "NOT FOUND #14: formato de arquivo.[13] ABREU" }
$found15 = $d.Content.Find.Execute("Extensionista:- Disseminar a cultura", $true, $false, $false, $false, $false, $true, 1, $false, "Extensionista:^l- Disseminar a cultura", 2)
if (-not $found15) { Write-Host "NOT FOUND #15: Extensionista:- Disseminar a cultura" }
$found16 = $d.Content.Find.Execute("sociedade impactada.- Estimular os estudantes", $true, $false, $false, $false, $false, $true, 1, $false, "sociedade impactada.^l- Estimular os estudantes", 2)
if (-not $found16) { Write-Host "NOT FOUND #16: sociedade impactada.- Estimular os estudantes" }
$found17 = $d.Content.Find.Execute("esfera “Soft Skills”.A atividade consiste", $true, $false, $false, $false, $false, $true, 1, $false, "esfera “Soft Skills”.^lA atividade consiste", 2)
if (-not $found17) { Write-Host "NOT FOUND #17: esfera “Soft Skills”.A atividade consiste" }
$found18 = $d.Content.Find.Execute("profissionais etc..).Etapas:1.Planejamento da Oficina", $true, $false, $false, $false, $false, $true, 1, $false, "profissionais etc..).^l^lEtapas:^l1.Planejamento da Oficina", 2)
if (-not $found18) { Write-Host "NOT FOUND #18: profissionais etc..).Etapas:1.Planejamento da Oficina" }
$found19 = $d.Content.Find.Execute("professor da disciplina.2.Realização da Oficina", $true, $false, $false, $false, $false, $true, 1, $false, "professor da disciplina.^l2.Realização da Oficina", 2)
if (-not $found19) { Write-Host "NOT FOUND #19: professor da disciplina.2.Realização da Oficina" }
$found20 = $d.Content.Find.Execute("totais de aplicação.3.Preparação de Relatos", $true, $false, $false, $false, $false, $true, 1, $false, "totais de aplicação.^l3.Preparação de Relatos", 2)
if (-not $found20) { Write-Host "NOT FOUND #20: totais de aplicação.3.Preparação de Relatos" }
$found21 = $d.Content.Find.Execute("para a comunidade.4.Autoavaliação pelo Grupo", $true, $false, $false, $false, $false, $true, 1, $false, "para a comunidade.^l4.Autoavaliação pelo Grupo", 2)
if (-not $found21) { Write-Host "NOT FOUND #21: para a comunidade.4.Autoavaliação pelo Grupo" }

Write-Host "Done"
